# "Antes de aplicar el crossover" -> "despues": refresh the DE trial data
# on both sheets and tidy up the view state to match the author's session.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: the "your algorithm name" row becomes the "DE" row, and its
# F1..F10 statistics (C4:L4) are filled in from the Hoja1 trial means.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("B4").Value = "DE"

$ws1.Cells.Item(4, "C").Value = 80501.2695
$ws1.Cells.Item(4, "D").Value = 958286.595
$ws1.Cells.Item(4, "E").Value = 1493601325
$ws1.Cells.Item(4, "F").Value = 114310.7055
$ws1.Cells.Item(4, "G").Value = 38134.527
$ws1.Cells.Item(4, "H").Value = 51408609500
$ws1.Cells.Item(4, "I").Value = 3310.82555
$ws1.Cells.Item(4, "J").Value = 21.1287395
$ws1.Cells.Item(4, "K").Value = 470.677785
$ws1.Cells.Item(4, "L").Value = 758.308475

# Column C did not have an explicit width before; the author set one.
$ws1.Columns.Item(3).ColumnWidth = 8

$ws1.Range("L5").Select()

# ---------------------------------------------------------------------
# Hoja1: a brand-new batch of 20 algorithm runs, ten metrics (A:J) wide.
# Columns D:J did not exist yet for this batch, so give them the same
# scientific-notation display the existing A:C columns already use.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Hoja1")

$ws2.Range("D2:J21").NumberFormat = "0.00E+00"

# row 2
$ws2.Cells.Item(2, "A").Value = 67923210000
$ws2.Cells.Item(2, "B").Value = 105783300000
$ws2.Cells.Item(2, "C").Value = 1781265000000000
$ws2.Cells.Item(2, "D").Value = 135666200000
$ws2.Cells.Item(2, "E").Value = 38593470000
$ws2.Cells.Item(2, "F").Value = 44890440000000000
$ws2.Cells.Item(2, "G").Value = 3407259000
$ws2.Cells.Item(2, "H").Value = 21125470
$ws2.Cells.Item(2, "I").Value = 386125800
$ws2.Cells.Item(2, "J").Value = 687723200

# row 3
$ws2.Cells.Item(3, "A").Value = 79170560000
$ws2.Cells.Item(3, "B").Value = 97323860000
$ws2.Cells.Item(3, "C").Value = 1573398000000000
$ws2.Cells.Item(3, "D").Value = 136289600000
$ws2.Cells.Item(3, "E").Value = 28633720000
$ws2.Cells.Item(3, "F").Value = 66996830000000000
$ws2.Cells.Item(3, "G").Value = 3507623000
$ws2.Cells.Item(3, "H").Value = 21229310
$ws2.Cells.Item(3, "I").Value = 487545700
$ws2.Cells.Item(3, "J").Value = 770658000

# row 4
$ws2.Cells.Item(4, "A").Value = 74183740000
$ws2.Cells.Item(4, "B").Value = 90045750000
$ws2.Cells.Item(4, "C").Value = 1068987000000000
$ws2.Cells.Item(4, "D").Value = 145475500000
$ws2.Cells.Item(4, "E").Value = 39884030000
$ws2.Cells.Item(4, "F").Value = 45970230000000000
$ws2.Cells.Item(4, "G").Value = 3191990000
$ws2.Cells.Item(4, "H").Value = 20977280
$ws2.Cells.Item(4, "I").Value = 504693300
$ws2.Cells.Item(4, "J").Value = 730727600

# row 5
$ws2.Cells.Item(5, "A").Value = 90705870000
$ws2.Cells.Item(5, "B").Value = 107190900000
$ws2.Cells.Item(5, "C").Value = 1647230000000000
$ws2.Cells.Item(5, "D").Value = 117950700000
$ws2.Cells.Item(5, "E").Value = 37283010000
$ws2.Cells.Item(5, "F").Value = 62494430000000000
$ws2.Cells.Item(5, "G").Value = 3137175000
$ws2.Cells.Item(5, "H").Value = 21145180
$ws2.Cells.Item(5, "I").Value = 513645400
$ws2.Cells.Item(5, "J").Value = 785434100

# row 6
$ws2.Cells.Item(6, "A").Value = 86409130000
$ws2.Cells.Item(6, "B").Value = 103934300000
$ws2.Cells.Item(6, "C").Value = 1416933000000000
$ws2.Cells.Item(6, "D").Value = 95979790000
$ws2.Cells.Item(6, "E").Value = 42329010000
$ws2.Cells.Item(6, "F").Value = 61236100000000000
$ws2.Cells.Item(6, "G").Value = 2779774000
$ws2.Cells.Item(6, "H").Value = 21157070
$ws2.Cells.Item(6, "I").Value = 488941000
$ws2.Cells.Item(6, "J").Value = 883090600

# row 7
$ws2.Cells.Item(7, "A").Value = 78785120000
$ws2.Cells.Item(7, "B").Value = 94051690000
$ws2.Cells.Item(7, "C").Value = 1688509000000000
$ws2.Cells.Item(7, "D").Value = 131347000000
$ws2.Cells.Item(7, "E").Value = 33219950000
$ws2.Cells.Item(7, "F").Value = 63260470000000000
$ws2.Cells.Item(7, "G").Value = 3771613000
$ws2.Cells.Item(7, "H").Value = 21153610
$ws2.Cells.Item(7, "I").Value = 438520600
$ws2.Cells.Item(7, "J").Value = 732375300

# row 8
$ws2.Cells.Item(8, "A").Value = 68678010000
$ws2.Cells.Item(8, "B").Value = 67942280000
$ws2.Cells.Item(8, "C").Value = 1684294000000000
$ws2.Cells.Item(8, "D").Value = 75206040000
$ws2.Cells.Item(8, "E").Value = 39173190000
$ws2.Cells.Item(8, "F").Value = 29199280000000000
$ws2.Cells.Item(8, "G").Value = 3379354000
$ws2.Cells.Item(8, "H").Value = 20999680
$ws2.Cells.Item(8, "I").Value = 518372700
$ws2.Cells.Item(8, "J").Value = 716323600

# row 9
$ws2.Cells.Item(9, "A").Value = 79523310000
$ws2.Cells.Item(9, "B").Value = 113467300000
$ws2.Cells.Item(9, "C").Value = 1315186000000000
$ws2.Cells.Item(9, "D").Value = 117517600000
$ws2.Cells.Item(9, "E").Value = 42752220000
$ws2.Cells.Item(9, "F").Value = 82330690000000000
$ws2.Cells.Item(9, "G").Value = 3545326000
$ws2.Cells.Item(9, "H").Value = 21239680
$ws2.Cells.Item(9, "I").Value = 463204400
$ws2.Cells.Item(9, "J").Value = 804896400

# row 10
$ws2.Cells.Item(10, "A").Value = 72482060000
$ws2.Cells.Item(10, "B").Value = 66781060000
$ws2.Cells.Item(10, "C").Value = 1885020000000000
$ws2.Cells.Item(10, "D").Value = 131450700000
$ws2.Cells.Item(10, "E").Value = 38417300000
$ws2.Cells.Item(10, "F").Value = 46890850000000000
$ws2.Cells.Item(10, "G").Value = 3039341000
$ws2.Cells.Item(10, "H").Value = 21123310
$ws2.Cells.Item(10, "I").Value = 434853200
$ws2.Cells.Item(10, "J").Value = 621198100

# row 11
$ws2.Cells.Item(11, "A").Value = 80314760000
$ws2.Cells.Item(11, "B").Value = 90830880000
$ws2.Cells.Item(11, "C").Value = 989423700000000
$ws2.Cells.Item(11, "D").Value = 100053700000
$ws2.Cells.Item(11, "E").Value = 41646840000
$ws2.Cells.Item(11, "F").Value = 42012910000000000
$ws2.Cells.Item(11, "G").Value = 3442875000
$ws2.Cells.Item(11, "H").Value = 21156360
$ws2.Cells.Item(11, "I").Value = 449134900
$ws2.Cells.Item(11, "J").Value = 831096700

# row 12
$ws2.Cells.Item(12, "A").Value = 56426660000
$ws2.Cells.Item(12, "B").Value = 90897390000
$ws2.Cells.Item(12, "C").Value = 1592365000000000
$ws2.Cells.Item(12, "D").Value = 98572030000
$ws2.Cells.Item(12, "E").Value = 37295120000
$ws2.Cells.Item(12, "F").Value = 38227610000000000
$ws2.Cells.Item(12, "G").Value = 3440532000
$ws2.Cells.Item(12, "H").Value = 21165280
$ws2.Cells.Item(12, "I").Value = 477545100
$ws2.Cells.Item(12, "J").Value = 809970900

# row 13
$ws2.Cells.Item(13, "A").Value = 95982760000
$ws2.Cells.Item(13, "B").Value = 109057400000
$ws2.Cells.Item(13, "C").Value = 1136752000000000
$ws2.Cells.Item(13, "D").Value = 75314290000
$ws2.Cells.Item(13, "E").Value = 38856520000
$ws2.Cells.Item(13, "F").Value = 29390140000000000
$ws2.Cells.Item(13, "G").Value = 2629716000
$ws2.Cells.Item(13, "H").Value = 21067850
$ws2.Cells.Item(13, "I").Value = 424454200
$ws2.Cells.Item(13, "J").Value = 752204300

# row 14
$ws2.Cells.Item(14, "A").Value = 75021880000
$ws2.Cells.Item(14, "B").Value = 92215610000
$ws2.Cells.Item(14, "C").Value = 1922128000000000
$ws2.Cells.Item(14, "D").Value = 124575500000
$ws2.Cells.Item(14, "E").Value = 35834210000
$ws2.Cells.Item(14, "F").Value = 43075980000000000
$ws2.Cells.Item(14, "G").Value = 2958782000
$ws2.Cells.Item(14, "H").Value = 21146580
$ws2.Cells.Item(14, "I").Value = 473325000
$ws2.Cells.Item(14, "J").Value = 694924500

# row 15
$ws2.Cells.Item(15, "A").Value = 89418930000
$ws2.Cells.Item(15, "B").Value = 98386800000
$ws2.Cells.Item(15, "C").Value = 1994117000000000
$ws2.Cells.Item(15, "D").Value = 137707900000
$ws2.Cells.Item(15, "E").Value = 39300840000
$ws2.Cells.Item(15, "F").Value = 47876410000000000
$ws2.Cells.Item(15, "G").Value = 2850027000
$ws2.Cells.Item(15, "H").Value = 21158060
$ws2.Cells.Item(15, "I").Value = 499668200
$ws2.Cells.Item(15, "J").Value = 757281300

# row 16
$ws2.Cells.Item(16, "A").Value = 66729390000
$ws2.Cells.Item(16, "B").Value = 116480300000
$ws2.Cells.Item(16, "C").Value = 1632864000000000
$ws2.Cells.Item(16, "D").Value = 96434310000
$ws2.Cells.Item(16, "E").Value = 33035270000
$ws2.Cells.Item(16, "F").Value = 37698990000000000
$ws2.Cells.Item(16, "G").Value = 3345014000
$ws2.Cells.Item(16, "H").Value = 21102260
$ws2.Cells.Item(16, "I").Value = 500731100
$ws2.Cells.Item(16, "J").Value = 775654800

# row 17
$ws2.Cells.Item(17, "A").Value = 71083850000
$ws2.Cells.Item(17, "B").Value = 102142100000
$ws2.Cells.Item(17, "C").Value = 1721994000000000
$ws2.Cells.Item(17, "D").Value = 78681050000
$ws2.Cells.Item(17, "E").Value = 40688440000
$ws2.Cells.Item(17, "F").Value = 49647750000000000
$ws2.Cells.Item(17, "G").Value = 3683364000
$ws2.Cells.Item(17, "H").Value = 21087400
$ws2.Cells.Item(17, "I").Value = 487181700
$ws2.Cells.Item(17, "J").Value = 715069000

# row 18
$ws2.Cells.Item(18, "A").Value = 81350200000
$ws2.Cells.Item(18, "B").Value = 88269870000
$ws2.Cells.Item(18, "C").Value = 1415495000000000
$ws2.Cells.Item(18, "D").Value = 127932800000
$ws2.Cells.Item(18, "E").Value = 33757330000
$ws2.Cells.Item(18, "F").Value = 55834090000000000
$ws2.Cells.Item(18, "G").Value = 3429200000
$ws2.Cells.Item(18, "H").Value = 21178770
$ws2.Cells.Item(18, "I").Value = 485860400
$ws2.Cells.Item(18, "J").Value = 698167100

# row 19
$ws2.Cells.Item(19, "A").Value = 59368390000
$ws2.Cells.Item(19, "B").Value = 100386100000
$ws2.Cells.Item(19, "C").Value = 830615800000000
$ws2.Cells.Item(19, "D").Value = 123089500000
$ws2.Cells.Item(19, "E").Value = 43225190000
$ws2.Cells.Item(19, "F").Value = 64746460000000000
$ws2.Cells.Item(19, "G").Value = 3887318000
$ws2.Cells.Item(19, "H").Value = 21116450
$ws2.Cells.Item(19, "I").Value = 500688500
$ws2.Cells.Item(19, "J").Value = 784649500

# row 20
$ws2.Cells.Item(20, "A").Value = 67207100000
$ws2.Cells.Item(20, "B").Value = 95003160000
$ws2.Cells.Item(20, "C").Value = 1398926000000000
$ws2.Cells.Item(20, "D").Value = 111449900000
$ws2.Cells.Item(20, "E").Value = 43462710000
$ws2.Cells.Item(20, "F").Value = 49594440000000000
$ws2.Cells.Item(20, "G").Value = 3289419000
$ws2.Cells.Item(20, "H").Value = 21055210
$ws2.Cells.Item(20, "I").Value = 390388100
$ws2.Cells.Item(20, "J").Value = 825973800

# row 21
$ws2.Cells.Item(21, "A").Value = 66582900000
$ws2.Cells.Item(21, "B").Value = 86383140000
$ws2.Cells.Item(21, "C").Value = 1176524000000000
$ws2.Cells.Item(21, "D").Value = 125520000000
$ws2.Cells.Item(21, "E").Value = 35302170000
$ws2.Cells.Item(21, "F").Value = 66798090000000000
$ws2.Cells.Item(21, "G").Value = 3500809000
$ws2.Cells.Item(21, "H").Value = 21189980
$ws2.Cells.Item(21, "I").Value = 488676400
$ws2.Cells.Item(21, "J").Value = 788750700

$ws2.Range("A2:A21").Select()

# Extend the "average of the 20 runs" row into the six new metric columns,
# the same shared =SUM(col2:col21)/20 pattern the A:C columns already use.
$ws2.Range("D22:J22").Formula = "=SUM(D2:D21)/20"
$ws2.Range("D22:J22").NumberFormat = "0.00E+00"

Write-Output "edit applied"
